$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-empty conditional probabilities (and their
#     complements) for the last CPT table (rows 17-20) ---
$ws.Range("D17").Value = 0.02297050938
$ws.Range("E17").Formula = "=1-D17"

$ws.Range("D18").Value = 0.08860053619
$ws.Range("E18").Formula = "=1-D18"

$ws.Range("D19").Value = 0.030991957099999998
$ws.Range("E19").Formula = "=1-D19"

$ws.Range("D20").Value = 0.1294369973
$ws.Range("E20").Formula = "=1-D20"

# --- Box the whole CPT table (B16:E20) in a thin border, matching the
#     other CPT tables above it on the sheet ---
$rng = $ws.Range("B16:E20")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# E17 picked up its border formatting from the header row above (E16)
# rather than from the newly-boxed block, so match that explicitly.
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122)  # xlPasteFormats

# --- Restore the scroll position / selection recorded in the saved file ---
$ws.Range("A35").Select()

Write-Host "done"
